# Mẫu 26 - remove the leftover unresolved "vnpt.SiteAddress" placeholder
# run that followed "Địa chỉ: " in the company-info paragraph.
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# wdFindContinue = 1, wdReplaceOne = 2
$find.Execute("vnpt.SiteAddress", $true, $false, $false, $false, $false, `
               $true, 1, $false, "", 2)
